$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 181, shifting existing rows 181-200 down to 182-201
$ws.Rows.Item(181).Insert()

# Populate the new row 181 with the new weekly record
$ws.Cells.Item(181, 1).Value = 3
$ws.Cells.Item(181, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(181, 3).Value = "Coquimbo"
$ws.Cells.Item(181, 4).Value = 44918
$ws.Cells.Item(181, 4).NumberFormat = $ws.Cells.Item(182, 4).NumberFormat
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 100112052
$ws.Cells.Item(181, 7).Value = "Albahaca"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 135
$ws.Cells.Item(181, 11).Value = 5000
$ws.Cells.Item(181, 12).Value = 6000
$ws.Cells.Item(181, 13).Value = 5481
$ws.Cells.Item(181, 14).Value = "$/docena de matas"
$ws.Cells.Item(181, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(181, 16).Value = 914
$ws.Cells.Item(181, 17).Value = 6
$ws.Cells.Item(181, 18).Value = "Hortaliza"
